$d = $word.ActiveDocument

# =====================================================================
# Change 1: underline the whole "La primera vez que entramos..." item
# (paragraph mark + its single run both get <w:u w:val="single"/>)
# =====================================================================
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("La primera vez que entramos al sistema")) {
        $p.Range.Font.Underline = 1
        break
    }
}

# =====================================================================
# Change 2: in the "Como  administrador  podemos..." item, underline
#   "los usuarios estará desactivados por defecto hasta que los active
#    el administrador."
# while leaving the text before and after it (including the final
# " También se podrán crear a otros administradores...") un-underlined.
# Doing this as three separate Find/Replace passes reproduces the same
# run layout as the source edit.
# =====================================================================
$f1 = $d.Content
$f1.Find.ClearFormatting()
$f1.Find.Replacement.ClearFormatting()
$f1.Find.Replacement.Font.Underline = 1
$f1.Find.Execute("los usuarios estará", $true, $false, $false, $false, $false, $true, 1, $false, "los usuarios estará", 2) | Out-Null

$f2 = $d.Content
$f2.Find.ClearFormatting()
$f2.Find.Replacement.ClearFormatting()
$f2.Find.Replacement.Font.Underline = 1
$f2.Find.Execute("n desactivados por ", $true, $false, $false, $false, $false, $true, 1, $false, "n desactivados por ", 2) | Out-Null

$f3 = $d.Content
$f3.Find.ClearFormatting()
$f3.Find.Replacement.ClearFormatting()
$f3.Find.Replacement.Font.Underline = 1
$f3.Find.Execute("defecto hasta que los active el administrador.", $true, $false, $false, $false, $false, $true, 1, $false, "defecto hasta que los active el administrador.", 2) | Out-Null

# =====================================================================
# Change 3: insert a brand-new empty paragraph (identical formatting to
# the existing empty paragraph that follows "Todos estos valores...")
# right before the "Los usuarios verán, por lo tanto..." item.
# =====================================================================
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "" -and $p.Range.ParagraphFormat.LeftIndent -eq 54) {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text.StartsWith("Los usuarios verán, por lo tanto")) {
            $p.Range.InsertParagraphAfter() | Out-Null
            break
        }
    }
}

Write-Output "done"
